# Add 2022-Q3 data
# --------------------------------------------------------------------------
# 1) Insert a new worksheet named "2022-Q3" right after "总计" (i.e. as the
#    2nd tab), pushing all the quarterly sheets down by one position.
# 2) Populate the new sheet with the 2022-Q3 fund-holding data.
# 3) Insert a new row into "总计" summarizing the 2022-Q3 quarter and shift
#    the existing summary rows (and their running index in column A) down.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Step 2: populate the new "2022-Q3" sheet
# ---------------------------------------------------------------------
# Copy the header formatting (bold, centered, bordered) used by every other
# quarterly sheet from the "总计" sheet's own header row, then overwrite the
# values with the fund-table column headers.
$summarySheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3Sheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$dataRows = @(
    @(0, "002670", "万家沪深300指数增强A",     "20.85", "94.06", "1.97", "0.4107", 10),
    @(1, "002671", "万家沪深300指数增强C",     "10.38", "94.06", "1.97", "0.2045", 10),
    @(2, "005635", "博时量化多策略股票A",      "2.09",  "84.79", "2.43", "0.0508", 4),
    @(3, "515300", "嘉实沪深300红利低波动ETF", "0.94",  "99.19", "3.46", "0.0325", 4),
    @(4, "510290", "南方上证380ETF",           "1.55",  "99.03", "0.87", "0.0135", 9),
    @(5, "009107", "嘉合同顺智选股票C",        "0.30",  "85.51", "3.07", "0.0092", 7),
    @(6, "005636", "博时量化多策略股票C",      "0.11",  "84.79", "2.43", "0.0027", 4),
    @(7, "006992", "嘉合锦创优势精选混合",     "0.08",  "84.53", "3.24", "0.0026", 5),
    @(8, "009106", "嘉合同顺智选股票A",        "0.01",  "85.51", "3.07", "0.0003", 7)
)

$r = 2
foreach ($row in $dataRows) {
    # Column A: running index (plain number)
    $q3Sheet.Cells.Item($r, 1).Value = $row[0]

    # Columns B-G: text values (fund code / name / size / position / % / value)
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q3Sheet.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }

    # Column H: rank (plain number)
    $q3Sheet.Cells.Item($r, 8).Value = $row[7]

    $r++
}

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet with a new 2022-Q3 row
# ---------------------------------------------------------------------
# Insert a blank row at row 2 (pushes every existing quarter row down by one).
$summarySheet.Rows.Item(2).Insert()

# Re-number the running index (column A) of all the rows that got shifted
# down: Rows.Insert physically moved each row's existing cells (including
# column A's old running index) down by one without changing their values,
# so bump each of those carried-over index values by 1 in place.
$lastRow = 9
for ($row = 3; $row -le $lastRow; $row++) {
    $cell = $summarySheet.Cells.Item($row, 1)
    $oldIndex = $cell.Value()
    $cell.Value = $oldIndex + 1
}

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 9
$summarySheet.Cells.Item(2, 4).Value = 0.73
